# Generate Report for Handoff
# Replace the old GUID-based file identifiers with the new one, and refresh
# the handoff/target timestamps & generated xlf file names, across all three
# worksheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$oldGuid = "31807cc8-aaab-4b19-809b-4ea70eaf228f"
$newGuid = "24f3fbdb-43db-46bb-a4aa-3a64f07f679e"

$oldHash = "6af30091f61af7c11a4195c4fba221b33bf7f9af"
$newHash = "02db8b92cf30802664081aa8dbe6dc337d4cbd24"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"

foreach ($hl in $wsOverview.Hyperlinks) {
    $hl.TextToDisplay = "e2e\$newGuid.md"
}

$wsOverview.Range("G2").Value = "2016-09-05 21:07:34"

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "$newGuid.md"

foreach ($hl in $wsZh.Hyperlinks) {
    $hl.TextToDisplay = "$newGuid.md"
}

$wsZh.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-09-05 21:07:29"

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "$newGuid.md"

foreach ($hl in $wsDe.Hyperlinks) {
    $hl.TextToDisplay = "$newGuid.md"
}

$wsDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = "2016-09-05 21:07:34"
